# Updates cryptos list values (Price column D, Volume(1h) column E)
# Values are written as text (matching the inlineStr cells in the source)
# by briefly forcing a text number format, then clearing the format so the
# cell keeps no explicit style (matches the original, unstyled cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '57.933.74'
Set-TextValue 'E2' '  -3.84%  '
Set-TextValue 'D3' '2.292.60'
Set-TextValue 'E3' '  -4.89%  '
Set-TextValue 'E4' '  -0.02%  '
Set-TextValue 'D5' '538.87'
Set-TextValue 'E5' '  -3.62%  '
Set-TextValue 'D6' '131.04'
Set-TextValue 'E6' '  -3.42%  '
Set-TextValue 'E7' '  -0.02%  '
Set-TextValue 'D8' '0.568'
Set-TextValue 'E8' '  -3.54%  '
Set-TextValue 'D9' '2.291.50'
Set-TextValue 'E9' '  -4.83%  '
Set-TextValue 'E10' '  -5.43%  '
Set-TextValue 'D11' '5.47'
Set-TextValue 'E11' '  -2.42%  '
Set-TextValue 'E12' '  -0.27%  '
Set-TextValue 'D13' '0.331'
Set-TextValue 'E13' '  -5.26%  '
Set-TextValue 'D14' '23.57'
Set-TextValue 'E14' '  -4.54%  '
Set-TextValue 'D15' '2.700.30'
Set-TextValue 'E15' '  -4.82%  '
Set-TextValue 'D16' '57.880.13'
Set-TextValue 'E16' '  -3.79%  '
Set-TextValue 'E17' '  -4.63%  '
Set-TextValue 'D18' '2.270.32'
Set-TextValue 'E18' '  -4.39%  '
Set-TextValue 'D19' '10.60'
Set-TextValue 'E19' '  -5.49%  '
Set-TextValue 'D20' '4.25'
Set-TextValue 'E20' '  -6.63%  '
Set-TextValue 'D21' '313.73'
Set-TextValue 'E21' '  -3.50%  '
Set-TextValue 'D22' '6.43'
Set-TextValue 'E22' '  -5.80%  '
Set-TextValue 'E23' '  +0.03%  '
Set-TextValue 'D24' '63.04'
Set-TextValue 'E24' '  -2.43%  '
Set-TextValue 'D25' '0.167'
Set-TextValue 'E25' '  -4.39%  '
Set-TextValue 'D26' '0.999'
Set-TextValue 'E26' '  -0.11%  '
Set-TextValue 'D27' '7.99'
Set-TextValue 'E27' '  -6.32%  '
Set-TextValue 'D28' '1.30'
Set-TextValue 'E28' '  -6.01%  '
Set-TextValue 'D29' '169.76'
Set-TextValue 'E29' '  -0.52%  '
Set-TextValue 'D30' '1.72'
Set-TextValue 'E30' '  -4.74%  '
Set-TextValue 'D31' '0.0₃0721'
Set-TextValue 'E31' '  -6.74%  '
Set-TextValue 'E32' '  -0.35%  '
Set-TextValue 'D33' '5.72'
Set-TextValue 'E33' '  -6.59%  '
Set-TextValue 'D34' '0.380'
Set-TextValue 'E34' '  -5.83%  '
Set-TextValue 'E35' '  -0.01%  '
Set-TextValue 'D36' '17.70'
Set-TextValue 'E36' '  -3.63%  '
Set-TextValue 'D37' '1.00'
Set-TextValue 'E37' '  -0.01%  '
Set-TextValue 'D38' '1.24'
Set-TextValue 'E38' '  -7.65%  '
Set-TextValue 'D39' '3.90'
Set-TextValue 'E39' '  -6.45%  '
Set-TextValue 'D40' '37.85'
Set-TextValue 'E40' '  -2.05%  '
Set-TextValue 'D41' '1.49'
Set-TextValue 'E41' '  -6.84%  '
Set-TextValue 'D42' '141.90'
Set-TextValue 'E42' '  -4.18%  '
Set-TextValue 'D43' '289.96'
Set-TextValue 'E43' '  -10.52%  '
Set-TextValue 'D44' '3.40'
Set-TextValue 'E44' '  -4.59%  '
Set-TextValue 'D45' '0.0944'
Set-TextValue 'E45' '  -2.65%  '
Set-TextValue 'D46' '0.0498'
Set-TextValue 'E46' '  -3.46%  '
Set-TextValue 'D47' '0.554'
Set-TextValue 'E47' '  -3.57%  '
Set-TextValue 'D48' '18.25'
Set-TextValue 'E48' '  -8.27%  '
Set-TextValue 'D49' '0.0211'
Set-TextValue 'E49' '  -4.81%  '
Set-TextValue 'D50' '16.56'
Set-TextValue 'E50' '  -3.72%  '
Set-TextValue 'E51' '  -0.95%  '
